# Commit: "updated few testcases to yes"
# Sets the Execution_Status column (F) on the TestSuite sheet to "No"
# for rows 26-45 (rows 2-25 stay "Yes").
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestSuite")
$ws.Activate()

foreach ($r in 26..45) {
    $ws.Range("F$r").Value = "No"
}

# Match the author's final cursor position/selection in the saved file.
$ws.Range("F48").Select()
